$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Birth Date) holds free-text date strings that are NOT valid
# ISO dates, formatted as DD/MM/YYYY. Force the column to Text first so
# Excel does not auto-convert the literal strings into date serials.
$ws.Range("C2:C31").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Wayne Lopez"
$ws.Range("B2").Value = "9**********4"
$ws.Range("C2").Value = "12/11/1992"
$ws.Range("D2").Value = "USCGC Willis`nFPO AP 35902"
$ws.Range("E2").Value = "***-*****89"
$ws.Range("H2").Value = "Smithburgh"
$ws.Range("J2").Value = 52
$ws.Range("M2").Value = 3582
$ws.Range("O2").Value = "population"

# Row 3
$ws.Range("A3").Value = "Nichole Parsons"
$ws.Range("B3").Value = "9**********5"
$ws.Range("C3").Value = "30/05/1957"
$ws.Range("D3").Value = "97135 Boyd Glen`nMcdanielmouth, VA 10567"
$ws.Range("E3").Value = "***-*****90"
$ws.Range("H3").Value = "West Crystalfurt"
$ws.Range("J3").Value = 20
$ws.Range("M3").Value = 7177
$ws.Range("O3").Value = "cause"

# Row 4
$ws.Range("A4").Value = "Linda Bailey"
$ws.Range("B4").Value = "9**********6"
$ws.Range("C4").Value = "10/08/1997"
$ws.Range("D4").Value = "937 Foster Walks`nGrossstad, DE 27171"
$ws.Range("E4").Value = "***-*****01"
$ws.Range("H4").Value = "West Trevor"
$ws.Range("J4").Value = 65
$ws.Range("M4").Value = 9711
$ws.Range("O4").Value = "enjoy"

# Row 5
$ws.Range("A5").Value = "Joshua Flores"
$ws.Range("B5").Value = "9**********7"
$ws.Range("C5").Value = "17/02/1927"
$ws.Range("D5").Value = "PSC 4065, Box 0678`nAPO AE 60035"
$ws.Range("E5").Value = "***-*****12"
$ws.Range("H5").Value = "West Bethbury"
$ws.Range("J5").Value = 20
$ws.Range("M5").Value = 4186
$ws.Range("O5").Value = "record"

# Row 6
$ws.Range("A6").Value = "Robert Salas"
$ws.Range("B6").Value = "9**********8"
$ws.Range("C6").Value = "28/10/1996"
$ws.Range("D6").Value = "2001 Nicole Avenue`nPatrickfurt, PW 94092"
$ws.Range("E6").Value = "***-*****23"
$ws.Range("H6").Value = "Lake Kennethstad"
$ws.Range("J6").Value = 29
$ws.Range("M6").Value = 6285
$ws.Range("O6").Value = "growth"

# Row 7
$ws.Range("A7").Value = "Jason Johnson"
$ws.Range("B7").Value = "9**********9"
$ws.Range("C7").Value = "24/12/1992"
$ws.Range("D7").Value = "696 Thomas Views Suite 830`nEast Marymouth, WI 73270"
$ws.Range("E7").Value = "***-*****34"
$ws.Range("H7").Value = "Jacksonmouth"
$ws.Range("J7").Value = 90
$ws.Range("M7").Value = 6755
$ws.Range("O7").Value = "wear"

# Row 8
$ws.Range("A8").Value = "Alejandro Francis"
$ws.Range("B8").Value = "9**********0"
$ws.Range("C8").Value = "04/09/1962"
$ws.Range("D8").Value = "754 Perez Prairie`nNorth Rogerfort, KY 48698"
$ws.Range("E8").Value = "***-*****45"
$ws.Range("H8").Value = "Powerstown"
$ws.Range("J8").Value = 76
$ws.Range("M8").Value = 7073
$ws.Range("O8").Value = "analysis"

# Row 9
$ws.Range("A9").Value = "Stephanie Roberts"
$ws.Range("B9").Value = "9**********1"
$ws.Range("C9").Value = "19/08/1930"
$ws.Range("D9").Value = "494 Hess Place`nFosterburgh, RI 34180"
$ws.Range("E9").Value = "***-*****56"
$ws.Range("H9").Value = "Aprilburgh"
$ws.Range("J9").Value = 47
$ws.Range("M9").Value = 2351
$ws.Range("O9").Value = "information"

# Row 10
$ws.Range("A10").Value = "Karen Burton"
$ws.Range("B10").Value = "0**********2"
$ws.Range("C10").Value = "10/12/1934"
$ws.Range("D10").Value = "69310 Jessica Rest`nMichaelstad, SC 81770"
$ws.Range("E10").Value = "***-*****78"
$ws.Range("H10").Value = "Lake Roberthaven"
$ws.Range("J10").Value = 89
$ws.Range("M10").Value = 9135
$ws.Range("O10").Value = "hit"

# Row 11
$ws.Range("A11").Value = "Justin Rodriguez"
$ws.Range("B11").Value = "0**********3"
$ws.Range("C11").Value = "24/06/1934"
$ws.Range("D11").Value = "Unit 0240 Box 1862`nDPO AP 15594"
$ws.Range("E11").Value = "***-*****89"
$ws.Range("H11").Value = "Georgeton"
$ws.Range("J11").Value = 74
$ws.Range("M11").Value = 9272
$ws.Range("O11").Value = "man"

# Row 12
$ws.Range("A12").Value = "Maureen Rivera"
$ws.Range("B12").Value = "0**********4"
$ws.Range("C12").Value = "10/04/1957"
$ws.Range("D12").Value = "9839 Mcmillan Lights Apt. 704`nJohnmouth, RI 68936"
$ws.Range("E12").Value = "***-*****01"
$ws.Range("H12").Value = "East Peter"
$ws.Range("J12").Value = 25
$ws.Range("M12").Value = 4605
$ws.Range("O12").Value = "election"

# Row 13
$ws.Range("A13").Value = "Darryl Weiss"
$ws.Range("B13").Value = "0**********5"
$ws.Range("C13").Value = "04/09/1930"
$ws.Range("D13").Value = "Unit 5002 Box 8334`nDPO AA 76006"
$ws.Range("E13").Value = "***-*****12"
$ws.Range("H13").Value = "West Molly"
$ws.Range("J13").Value = 95
$ws.Range("M13").Value = 8922
$ws.Range("O13").Value = "help"

# Row 14
$ws.Range("A14").Value = "Daniel Walker"
$ws.Range("B14").Value = "0**********6"
$ws.Range("C14").Value = "10/12/1929"
$ws.Range("D14").Value = "935 Morgan Drives`nRyanport, OR 60678"
$ws.Range("E14").Value = "***-*****23"
$ws.Range("H14").Value = "Anthonytown"
$ws.Range("J14").Value = 25
$ws.Range("M14").Value = 7886
$ws.Range("O14").Value = "difference"

# Row 15
$ws.Range("A15").Value = "Stephen Clark"
$ws.Range("B15").Value = "0**********7"
$ws.Range("C15").Value = "25/08/1985"
$ws.Range("D15").Value = "21848 Avila Court`nWest Brandon, WA 47870"
$ws.Range("E15").Value = "***-*****34"
$ws.Range("H15").Value = "Caitlinmouth"
$ws.Range("J15").Value = 34
$ws.Range("M15").Value = 7500
$ws.Range("O15").Value = "pattern"

# Row 16
$ws.Range("A16").Value = "Travis Mccall"
$ws.Range("B16").Value = "0**********8"
$ws.Range("C16").Value = "14/04/1951"
$ws.Range("D16").Value = "3511 Thomas Way`nPort Jamesstad, NM 25197"
$ws.Range("E16").Value = "***-*****45"
$ws.Range("H16").Value = "East Samuel"
$ws.Range("J16").Value = 69
$ws.Range("M16").Value = 9736
$ws.Range("O16").Value = "I"

# Row 17
$ws.Range("A17").Value = "Joshua Garcia"
$ws.Range("B17").Value = "0**********9"
$ws.Range("C17").Value = "16/07/2006"
$ws.Range("D17").Value = "3293 Gary Drive Suite 787`nSouth Krista, NM 37393"
$ws.Range("E17").Value = "***-*****56"
$ws.Range("H17").Value = "West Mariehaven"
$ws.Range("J17").Value = 54
$ws.Range("M17").Value = 4714
$ws.Range("O17").Value = "collection"

# Row 18
$ws.Range("A18").Value = "Christopher Carrillo"
$ws.Range("B18").Value = "0**********0"
$ws.Range("C18").Value = "13/02/1983"
$ws.Range("D18").Value = "8515 Thomas Parks`nClineborough, NV 61599"
$ws.Range("E18").Value = "***-*****67"
$ws.Range("H18").Value = "Port Cindy"
$ws.Range("J18").Value = 25
$ws.Range("M18").Value = 2942
$ws.Range("O18").Value = "outside"

# Row 19
$ws.Range("A19").Value = "Lindsey Klein"
$ws.Range("B19").Value = "0**********1"
$ws.Range("C19").Value = "09/05/1980"
$ws.Range("D19").Value = "890 John Corners Apt. 904`nLake Lisaland, NY 94693"
$ws.Range("E19").Value = "***-*****78"
$ws.Range("H19").Value = "Danielland"
$ws.Range("J19").Value = 18
$ws.Range("M19").Value = 8815
$ws.Range("O19").Value = "example"

# Row 20
$ws.Range("A20").Value = "Ariel Rubio"
$ws.Range("B20").Value = "1**********2"
$ws.Range("C20").Value = "18/07/1980"
$ws.Range("D20").Value = "8090 Jason Rue Suite 045`nLake Stevenmouth, AS 90385"
$ws.Range("E20").Value = "***-*****89"
$ws.Range("H20").Value = "South Kimberlyfurt"
$ws.Range("J20").Value = 38
$ws.Range("M20").Value = 9754
$ws.Range("O20").Value = "no"

# Row 21
$ws.Range("A21").Value = "Alexis Hoover DVM"
$ws.Range("B21").Value = "1**********3"
$ws.Range("C21").Value = "07/03/1993"
$ws.Range("D21").Value = "20036 Lori Islands Suite 642`nLake Christopher, MO 43411"
$ws.Range("E21").Value = "***-*****90"
$ws.Range("H21").Value = "West Steventown"
$ws.Range("J21").Value = 56
$ws.Range("M21").Value = 6085
$ws.Range("O21").Value = "manager"

# Row 22
$ws.Range("A22").Value = "Jermaine Sims"
$ws.Range("B22").Value = "1**********4"
$ws.Range("C22").Value = "05/09/1941"
$ws.Range("D22").Value = "138 Roberts Ports`nPort Suzanne, PA 03681"
$ws.Range("E22").Value = "***-*****01"
$ws.Range("H22").Value = "East Jay"
$ws.Range("J22").Value = 93
$ws.Range("M22").Value = 2072
$ws.Range("O22").Value = "head"

# Row 23
$ws.Range("A23").Value = "Laurie Bell DVM"
$ws.Range("B23").Value = "1**********5"
$ws.Range("C23").Value = "23/03/1941"
$ws.Range("D23").Value = "Unit 2160 Box 6707`nDPO AE 26701"
$ws.Range("E23").Value = "***-*****12"
$ws.Range("H23").Value = "Stephaniestad"
$ws.Range("J23").Value = 79
$ws.Range("M23").Value = 8977
$ws.Range("O23").Value = "last"

# Row 24
$ws.Range("A24").Value = "Robert Duncan"
$ws.Range("B24").Value = "1**********6"
$ws.Range("C24").Value = "28/04/1995"
$ws.Range("D24").Value = "4296 King Creek Apt. 808`nNorth Ellen, AK 29042"
$ws.Range("E24").Value = "***-*****23"
$ws.Range("H24").Value = "Port Lacey"
$ws.Range("J24").Value = 41
$ws.Range("M24").Value = 8937
$ws.Range("O24").Value = "lay"

# Row 25
$ws.Range("A25").Value = "Dawn Pham"
$ws.Range("B25").Value = "1**********7"
$ws.Range("C25").Value = "13/09/1971"
$ws.Range("D25").Value = "1533 Laurie Corners Apt. 680`nSouth Christopher, MD 13700"
$ws.Range("E25").Value = "***-*****34"
$ws.Range("H25").Value = "East Ronaldland"
$ws.Range("J25").Value = 91
$ws.Range("M25").Value = 7945
$ws.Range("O25").Value = "modern"

# Row 26
$ws.Range("A26").Value = "Kathy Powell"
$ws.Range("B26").Value = "1**********8"
$ws.Range("C26").Value = "11/03/1931"
$ws.Range("D26").Value = "586 Long Trail`nPort Laurieland, MA 84267"
$ws.Range("E26").Value = "***-*****45"
$ws.Range("H26").Value = "Lake Brian"
$ws.Range("J26").Value = 30
$ws.Range("M26").Value = 3367
$ws.Range("O26").Value = "world"

# Row 27
$ws.Range("A27").Value = "Kimberly Gray"
$ws.Range("B27").Value = "1**********9"
$ws.Range("C27").Value = "09/03/1962"
$ws.Range("D27").Value = "168 Philip Knolls`nKathrynside, ID 60944"
$ws.Range("E27").Value = "***-*****56"
$ws.Range("H27").Value = "Espinozaside"
$ws.Range("J27").Value = 95
$ws.Range("M27").Value = 5531
$ws.Range("O27").Value = "marriage"

# Row 28
$ws.Range("A28").Value = "Kathryn White"
$ws.Range("B28").Value = "1**********0"
$ws.Range("C28").Value = "21/01/1939"
$ws.Range("D28").Value = "729 Derek Highway Suite 964`nJeffreyhaven, CO 81240"
$ws.Range("E28").Value = "***-*****67"
$ws.Range("H28").Value = "Lake Bryan"
$ws.Range("J28").Value = 49
$ws.Range("M28").Value = 6700
$ws.Range("O28").Value = "decide"

# Row 29
$ws.Range("A29").Value = "Nicole Kidd"
$ws.Range("B29").Value = "1**********1"
$ws.Range("C29").Value = "20/10/1974"
$ws.Range("D29").Value = "0339 Cameron Brooks Apt. 639`nEast Chad, AS 35485"
$ws.Range("E29").Value = "***-*****78"
$ws.Range("H29").Value = "North Davidfurt"
$ws.Range("J29").Value = 77
$ws.Range("M29").Value = 4686
$ws.Range("O29").Value = "voice"

# Row 30
$ws.Range("A30").Value = "Justin Oneal"
$ws.Range("B30").Value = "2**********2"
$ws.Range("C30").Value = "25/10/1959"
$ws.Range("D30").Value = "892 Allen Terrace Suite 908`nAnnaville, KY 05160"
$ws.Range("E30").Value = "***-*****89"
$ws.Range("H30").Value = "Lake Jessebury"
$ws.Range("J30").Value = 35
$ws.Range("M30").Value = 7574
$ws.Range("O30").Value = "fish"

# Row 31
$ws.Range("A31").Value = "Janice Ibarra"
$ws.Range("B31").Value = "2**********3"
$ws.Range("C31").Value = "28/01/2000"
$ws.Range("D31").Value = "PSC 1569, Box 9608`nAPO AP 45134"
$ws.Range("E31").Value = "***-*****90"
$ws.Range("H31").Value = "Leemouth"
$ws.Range("J31").Value = 74
$ws.Range("M31").Value = 9315
$ws.Range("O31").Value = "hospital"

